# ============================================================
# Edit script: restructure PlayerPerformance workbook
#  - add 'Player Info' sheet before 'ODI Batting'
#  - add 'ODI Batting Extra' sheet after 'ODI Batting'
#  - rename ODI Batting!D1 MATCH_CARD_LINK -> MATCH_CODE and
#    replace the scorecard URLs with bare match codes
#  - drop the stray empty INNING_NUMBER cells (rows 50/65/87)
# ============================================================

$wb = $excel.ActiveWorkbook
$wsBatting = $wb.Worksheets.Item("ODI Batting")

# --- 1. New sheet: Player Info (placed before ODI Batting) ---
$wsInfo = $wb.Worksheets.Add($wsBatting)
$wsInfo.Name = "Player Info"

$infoHeader = $wsInfo.Range("A1:D1")
$infoHeaderVals = New-Object 'object[,]' 1,4
$infoHeaderVals[0,0] = "ID"
$infoHeaderVals[0,1] = "NAME"
$infoHeaderVals[0,2] = "BATTING_HAND"
$infoHeaderVals[0,3] = "BOWL_STYLE"
$infoHeader.Value = $infoHeaderVals
$infoHeader.Font.Bold = $true
$infoHeader.Borders.LineStyle = 1
$infoHeader.HorizontalAlignment = -4108
$infoHeader.VerticalAlignment = -4160

$infoRow2 = $wsInfo.Range("A2:D2")
$infoRow2.NumberFormat = "@"
$infoRow2Vals = New-Object 'object[,]' 1,4
$infoRow2Vals[0,0] = "3889"
$infoRow2Vals[0,1] = "Ajinkya M Rahane"
$infoRow2Vals[0,2] = "Right Handed"
$infoRow2Vals[0,3] = "Right Arm Medium"
$infoRow2.Value = $infoRow2Vals

$wsInfo.Range("A1").Select()

# --- 2. New sheet: ODI Batting Extra (placed after ODI Batting) ---
# Re-resolve the ODI Batting reference: the handle captured before the
# first Worksheets.Add() tracks a stale position once the sheet order shifts.
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsExtra = $wb.Worksheets.Add($null, $wsBatting)
$wsExtra.Name = "ODI Batting Extra"

$extraHeader = $wsExtra.Range("A1:F1")
$extraHeaderVals = New-Object 'object[,]' 1,6
$extraHeaderVals[0,0] = "MATCH_CODE"
$extraHeaderVals[0,1] = "BATTING_POSITION"
$extraHeaderVals[0,2] = "NUM_4"
$extraHeaderVals[0,3] = "NUM_6"
$extraHeaderVals[0,4] = "PERCENT_RUNS_OF_TOTAL"
$extraHeaderVals[0,5] = "MAN_OF_MATCH"
$extraHeader.Value = $extraHeaderVals
$extraHeader.Font.Bold = $true
$extraHeader.Borders.LineStyle = 1
$extraHeader.HorizontalAlignment = -4108
$extraHeader.VerticalAlignment = -4160

# Column A (MATCH_CODE) and F (MAN_OF_MATCH) are always populated as text
$extraA = $wsExtra.Range("A2:A21")
$extraA.NumberFormat = "@"
$extraAVals = New-Object 'object[,]' 20,1
$extraAVals[0,0] = "3954"
$extraAVals[1,0] = "3955"
$extraAVals[2,0] = "3978"
$extraAVals[3,0] = "4051"
$extraAVals[4,0] = "4052"
$extraAVals[5,0] = "4053"
$extraAVals[6,0] = "4056"
$extraAVals[7,0] = "4057"
$extraAVals[8,0] = "4065"
$extraAVals[9,0] = "4067"
$extraAVals[10,0] = "4069"
$extraAVals[11,0] = "4071"
$extraAVals[12,0] = "4074"
$extraAVals[13,0] = "4076"
$extraAVals[14,0] = "4126"
$extraAVals[15,0] = "4127"
$extraAVals[16,0] = "4128"
$extraAVals[17,0] = "4130"
$extraAVals[18,0] = "4133"
$extraAVals[19,0] = "4135"
$extraA.Value = $extraAVals

$extraF = $wsExtra.Range("F2:F21")
$extraFVals = New-Object 'object[,]' 20,1
$extraFVals[0,0] = "NO"
$extraFVals[1,0] = "NO"
$extraFVals[2,0] = "NO"
$extraFVals[3,0] = "NO"
$extraFVals[4,0] = "YES"
$extraFVals[5,0] = "NO"
$extraFVals[6,0] = "NO"
$extraFVals[7,0] = "NO"
$extraFVals[8,0] = "NO"
$extraFVals[9,0] = "NO"
$extraFVals[10,0] = "NO"
$extraFVals[11,0] = "NO"
$extraFVals[12,0] = "NO"
$extraFVals[13,0] = "NO"
$extraFVals[14,0] = "NO"
$extraFVals[15,0] = "NO"
$extraFVals[16,0] = "NO"
$extraFVals[17,0] = "NO"
$extraFVals[18,0] = "NO"
$extraFVals[19,0] = "NO"
$extraF.Value = $extraFVals

# Column B (BATTING_POSITION) holds a true number where known, blank otherwise
$wsExtra.Range("B2").Value = 1
$wsExtra.Range("B3").Value = 1
$wsExtra.Range("B6").Value = 1
$wsExtra.Range("B8").Value = 1
$wsExtra.Range("B9").Value = 1
$wsExtra.Range("B11").Value = 1
$wsExtra.Range("B12").Value = 1
$wsExtra.Range("B13").Value = 1
$wsExtra.Range("B15").Value = 1
$wsExtra.Range("B17").Value = 4
$wsExtra.Range("B18").Value = 4
$wsExtra.Range("B19").Value = 4
$wsExtra.Range("B20").Value = 4
$wsExtra.Range("B21").Value = 4

# Columns C, D (NUM_4 / NUM_6) and E (PERCENT_RUNS_OF_TOTAL) are text where known
$extraCDE = $wsExtra.Range("C2:E21")
$extraCDE.NumberFormat = "@"
$wsExtra.Range("C2").Value = "5"
$wsExtra.Range("D2").Value = "1"
$wsExtra.Range("E2").Value = "23.65%"
$wsExtra.Range("C3").Value = "3"
$wsExtra.Range("D3").Value = "0"
$wsExtra.Range("E3").Value = "7.43%"
$wsExtra.Range("C6").Value = "10"
$wsExtra.Range("D6").Value = "2"
$wsExtra.Range("E6").Value = "33.23%"
$wsExtra.Range("C8").Value = "7"
$wsExtra.Range("D8").Value = "0"
$wsExtra.Range("E8").Value = "33.71%"
$wsExtra.Range("C9").Value = "5"
$wsExtra.Range("D9").Value = "0"
$wsExtra.Range("E9").Value = "18.93%"
$wsExtra.Range("C11").Value = "0"
$wsExtra.Range("D11").Value = "0"
$wsExtra.Range("E11").Value = "1.78%"
$wsExtra.Range("C12").Value = "7"
$wsExtra.Range("D12").Value = "0"
$wsExtra.Range("E12").Value = "21.83%"
$wsExtra.Range("C13").Value = "9"
$wsExtra.Range("D13").Value = "0"
$wsExtra.Range("E13").Value = "23.81%"
$wsExtra.Range("C15").Value = "7"
$wsExtra.Range("D15").Value = "0"
$wsExtra.Range("E15").Value = "25.10%"
$wsExtra.Range("C18").Value = "0"
$wsExtra.Range("D18").Value = "0"
$wsExtra.Range("E18").Value = "3.63%"
$wsExtra.Range("C19").Value = "1"
$wsExtra.Range("D19").Value = "0"
$wsExtra.Range("E19").Value = "2.77%"
$wsExtra.Range("C20").Value = "0"
$wsExtra.Range("D20").Value = "0"
$wsExtra.Range("E20").Value = "2.92%"
$wsExtra.Range("C21").Value = "3"
$wsExtra.Range("D21").Value = "0"
$wsExtra.Range("E21").Value = "16.50%"

$wsExtra.Range("A1").Select()

# --- 3. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE ---
$wsBatting.Range("D1").Value = "MATCH_CODE"

$battingD = $wsBatting.Range("D2:D91")
$battingD.NumberFormat = "@"
$battingDVals = New-Object 'object[,]' 90,1
$battingDVals[0,0] = "3322"
$battingDVals[1,0] = "3323"
$battingDVals[2,0] = "3325"
$battingDVals[3,0] = "3327"
$battingDVals[4,0] = "3331"
$battingDVals[5,0] = "3335"
$battingDVals[6,0] = "3337"
$battingDVals[7,0] = "3341"
$battingDVals[8,0] = "3343"
$battingDVals[9,0] = "3346"
$battingDVals[10,0] = "3360"
$battingDVals[11,0] = "3435"
$battingDVals[12,0] = "3457"
$battingDVals[13,0] = "3459"
$battingDVals[14,0] = "3461"
$battingDVals[15,0] = "3463"
$battingDVals[16,0] = "3545"
$battingDVals[17,0] = "3587"
$battingDVals[18,0] = "3601"
$battingDVals[19,0] = "3603"
$battingDVals[20,0] = "3607"
$battingDVals[21,0] = "3609"
$battingDVals[22,0] = "3612"
$battingDVals[23,0] = "3619"
$battingDVals[24,0] = "3621"
$battingDVals[25,0] = "3624"
$battingDVals[26,0] = "3628"
$battingDVals[27,0] = "3643"
$battingDVals[28,0] = "3644"
$battingDVals[29,0] = "3645"
$battingDVals[30,0] = "3664"
$battingDVals[31,0] = "3667"
$battingDVals[32,0] = "3670"
$battingDVals[33,0] = "3672"
$battingDVals[34,0] = "3678"
$battingDVals[35,0] = "3680"
$battingDVals[36,0] = "3683"
$battingDVals[37,0] = "3688"
$battingDVals[38,0] = "3689"
$battingDVals[39,0] = "3692"
$battingDVals[40,0] = "3693"
$battingDVals[41,0] = "3696"
$battingDVals[42,0] = "3731"
$battingDVals[43,0] = "3735"
$battingDVals[44,0] = "3741"
$battingDVals[45,0] = "3744"
$battingDVals[46,0] = "3751"
$battingDVals[47,0] = "3760"
$battingDVals[48,0] = "3768"
$battingDVals[49,0] = "3775"
$battingDVals[50,0] = "3781"
$battingDVals[51,0] = "3786"
$battingDVals[52,0] = "3791"
$battingDVals[53,0] = "3795"
$battingDVals[54,0] = "3808"
$battingDVals[55,0] = "3812"
$battingDVals[56,0] = "3815"
$battingDVals[57,0] = "3817"
$battingDVals[58,0] = "3841"
$battingDVals[59,0] = "3844"
$battingDVals[60,0] = "3845"
$battingDVals[61,0] = "3848"
$battingDVals[62,0] = "3851"
$battingDVals[63,0] = "3874"
$battingDVals[64,0] = "3875"
$battingDVals[65,0] = "3876"
$battingDVals[66,0] = "3877"
$battingDVals[67,0] = "3951"
$battingDVals[68,0] = "3952"
$battingDVals[69,0] = "3953"
$battingDVals[70,0] = "3954"
$battingDVals[71,0] = "3955"
$battingDVals[72,0] = "3978"
$battingDVals[73,0] = "4051"
$battingDVals[74,0] = "4052"
$battingDVals[75,0] = "4053"
$battingDVals[76,0] = "4056"
$battingDVals[77,0] = "4057"
$battingDVals[78,0] = "4065"
$battingDVals[79,0] = "4067"
$battingDVals[80,0] = "4069"
$battingDVals[81,0] = "4071"
$battingDVals[82,0] = "4074"
$battingDVals[83,0] = "4076"
$battingDVals[84,0] = "4126"
$battingDVals[85,0] = "4127"
$battingDVals[86,0] = "4128"
$battingDVals[87,0] = "4130"
$battingDVals[88,0] = "4133"
$battingDVals[89,0] = "4135"
$battingD.Value = $battingDVals

# Drop the stray empty INNING_NUMBER placeholder cells
$wsBatting.Range("B50").ClearContents()
$wsBatting.Range("B65").ClearContents()
$wsBatting.Range("B87").ClearContents()

$wsBatting.Activate()
$wsBatting.Range("A1").Select()
